$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (Fitness) values for rows 2-252 change from 7573 to 7293
$ws.Range("C2:C252").Value = 7293
